# feat(maintenance): mark-as-serviced bulk update
#
# Simulates the "mark equipment as serviced" action: for each selected
# piece of equipment, the app writes the current timestamp (NOW()) into
# the "Дата последнего ТО" (last-maintenance date) column (I). The
# "Дата следующего ТО" column (J = I + interval) and the computed
# "Статус" column (K) recalculate automatically from the existing
# worksheet formulas.

$wb = $excel.ActiveWorkbook

# --- Sheet "ПК АСУ ТП" (PCs) ---------------------------------------------
$wsPc = $wb.Worksheets.Item("ПК АСУ ТП")

$pcServiced = @{
    7  = 46031.147372685184   # ККД-1 / ECS2261CLT03
    9  = 46031.151493055557   # ПНС-2 / ECS2261CLT05
    25 = 46031.14738425926    # ККД-1 / HiWatch (видеорегистратор)
    39 = 46031.144618055558   # ККД-2 / ECS5072CLT03
}

foreach ($row in $pcServiced.Keys) {
    $wsPc.Range("I$row").Value = $pcServiced[$row]
}

# --- Sheet "Шкафы АСУ ТП" (cabinets) -------------------------------------
$wsCab = $wb.Worksheets.Item("Шкафы АСУ ТП")

$cabServiced = @{
    260 = 46031.147418981483   # ККД-1 / 991CS110A01
    261 = 46031.144537037035   # ККД-2 / 997 Metso
    262 = 46031.150625000002   # Котельная / ЩУК-1 КВм-3.5 КБ
    263 = 46031.150659722225   # Котельная / ЩУК-2 КВм-3.5 КБ
    264 = 46031.150682870371   # Котельная / ЩУК-3 КВм-3.5 КБ
    265 = 46031.150717592594   # Котельная / ЩУК-4 КВм-3.5 КБ
    266 = 46031.150740740741   # Котельная / ЩУК-5 КВм-3.5 КБ
    267 = 46031.150775462964   # Котельная / ЩУК-6 КВм-3.5 КБ
    268 = 46031.15079861111    # Котельная / ЩУК-7 КВм-3.5 КБ
    269 = 46031.150833333333   # Котельная / ЩУК-8 А9.06.01.04
    270 = 46031.150856481479   # Котельная / ЩУК-9 А9.06.01.04
    271 = 46031.150891203702   # Котельная / ЩУК-10 А9.06.01.04
    272 = 46031.150914351849   # Котельная / ЩУК-11 А9.06.01.04
    273 = 46031.150949074072   # Котельная / ЩВО 1кВт
    274 = 46031.149988425925   # Котельная / СЩУ 301-02-АТХ.100
    277 = 46031.150011574071   # ПНС-2.1 / 996CS110A01
}

foreach ($row in $cabServiced.Keys) {
    $wsCab.Range("I$row").Value = $cabServiced[$row]
}
